$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data range entirely (A1:K5) so stale cells (D:K, rows 3-5) are removed.
$ws.Range("A1:K5").Clear()

# New header row
$ws.Range("A1").Value = "名称"
$ws.Range("B1").Value = "怪物"
$ws.Range("C1").Value = "等级"

# New data row (write B2 before A2 so the shared-string table order matches
# the target: qq人 must land before Hilichurl)
$ws.Range("B2").Value = "qq人"
$ws.Range("A2").Value = "Hilichurl"
$ws.Range("C2").Value = 90

# Update selection to match target
$ws.Range("C4").Select()
